# KHL stats refresh: append two new matches (Matches_SOG), roll forward the
# as_of_utc timestamp across the derived team-stat sheets, update the
# aggregate shot numbers for the four teams involved in the new matches,
# and bump Meta_ext's as_of_utc/build_version.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Matches_SOG: append rows 472 and 473
# ---------------------------------------------------------------------------
$matches = $wb.Worksheets.Item("Matches_SOG")

# uid column (A) must stay text, not auto-converted to a number.
$matches.Range("A472:A473").NumberFormat = "@"

$matches.Cells.Item(472, 1).Value = "897770"
$matches.Cells.Item(472, 2).Value = "2025-11-14T19:00:00"
$matches.Cells.Item(472, 3).Value = "Нефтехимик"
$matches.Cells.Item(472, 4).Value = "ЦСКА"
$matches.Cells.Item(472, 5).Value = 27
$matches.Cells.Item(472, 6).Value = 28
$matches.Cells.Item(472, 7).Value = "khl_text"

$matches.Cells.Item(473, 1).Value = "897769"
$matches.Cells.Item(473, 2).Value = "2025-11-14T19:30:00"
$matches.Cells.Item(473, 3).Value = "Спартак"
$matches.Cells.Item(473, 4).Value = "Северсталь"
$matches.Cells.Item(473, 5).Value = 27
$matches.Cells.Item(473, 6).Value = 25
$matches.Cells.Item(473, 7).Value = "khl_text"

# ---------------------------------------------------------------------------
# 2) Shots_HA: roll the as_of_utc column forward for every team (rows 2-23),
#    then update the on-goal totals/averages for the four teams that played.
# ---------------------------------------------------------------------------
$shotsHA = $wb.Worksheets.Item("Shots_HA")

for ($r = 2; $r -le 23; $r++) {
    $shotsHA.Cells.Item($r, 4).Value = "2025-11-14T19:30:00Z"
}

# Row 14: Нефтехимик
$shotsHA.Cells.Item(14, 5).Value = 26
$shotsHA.Cells.Item(14, 7).Value = 808
$shotsHA.Cells.Item(14, 8).Value = 886
$shotsHA.Cells.Item(14, 9).Value = 31.1
$shotsHA.Cells.Item(14, 10).Value = 34.1

# Row 17: Северсталь
$shotsHA.Cells.Item(17, 6).Value = 26
$shotsHA.Cells.Item(17, 11).Value = 854
$shotsHA.Cells.Item(17, 12).Value = 692
$shotsHA.Cells.Item(17, 13).Value = 32.8

# Row 19: Спартак
$shotsHA.Cells.Item(19, 5).Value = 27
$shotsHA.Cells.Item(19, 7).Value = 954
$shotsHA.Cells.Item(19, 8).Value = 746
$shotsHA.Cells.Item(19, 9).Value = 35.3
$shotsHA.Cells.Item(19, 10).Value = 27.6

# Row 23: ЦСКА
$shotsHA.Cells.Item(23, 6).Value = 22
$shotsHA.Cells.Item(23, 11).Value = 551
$shotsHA.Cells.Item(23, 12).Value = 624
$shotsHA.Cells.Item(23, 13).Value = 25

# ---------------------------------------------------------------------------
# 3) Shots_Summary: same as_of_utc roll-forward, then update team totals.
# ---------------------------------------------------------------------------
$shotsSummary = $wb.Worksheets.Item("Shots_Summary")

for ($r = 2; $r -le 23; $r++) {
    $shotsSummary.Cells.Item($r, 4).Value = "2025-11-14T19:30:00Z"
}

# Row 14: Нефтехимик
$shotsSummary.Cells.Item(14, 5).Value = 45
$shotsSummary.Cells.Item(14, 6).Value = 1325
$shotsSummary.Cells.Item(14, 7).Value = 1608
$shotsSummary.Cells.Item(14, 8).Value = 29.4
$shotsSummary.Cells.Item(14, 9).Value = 35.7

# Row 17: Северсталь
$shotsSummary.Cells.Item(17, 5).Value = 42
$shotsSummary.Cells.Item(17, 6).Value = 1323
$shotsSummary.Cells.Item(17, 7).Value = 1055
$shotsSummary.Cells.Item(17, 8).Value = 31.5

# Row 19: Спартак
$shotsSummary.Cells.Item(19, 5).Value = 42
$shotsSummary.Cells.Item(19, 6).Value = 1472
$shotsSummary.Cells.Item(19, 7).Value = 1280
$shotsSummary.Cells.Item(19, 8).Value = 35
$shotsSummary.Cells.Item(19, 9).Value = 30.5

# Row 23: ЦСКА
$shotsSummary.Cells.Item(23, 5).Value = 42
$shotsSummary.Cells.Item(23, 6).Value = 1021
$shotsSummary.Cells.Item(23, 7).Value = 1201
$shotsSummary.Cells.Item(23, 8).Value = 24.3

# ---------------------------------------------------------------------------
# 4) Meta_ext: bump as_of_utc and build_version.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Meta_ext")
$meta.Cells.Item(2, 2).Value = "2025-11-14T19:30:00Z"
$meta.Cells.Item(2, 4).Value = 66
